# finished the correction in columns id for xcel files sup34
#
# The header row (row 1) contains column-id labels that used a "vm..." naming
# prefix (a leftover/incorrect convention). This corrects every one of those
# 32 header labels (columns B:AG) to use the "spk..." prefix instead, e.g.
#   vmscpIsoStcDlat50            -> spkscpIsoStcDlat50
#   vmscpIsoStcDlat50Indisig     -> spkscpIsoStcDlat50Indisig
#   ...
# Only the literal text of the header cells changes; the row order of the
# columns, and all the numeric/string data in rows 2:21, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$renames = @{
    "B1"  = "spkscpIsoStcDlat50"
    "C1"  = "spkscpIsoStcDlat50Indisig"
    "D1"  = "spkscpIsoStcDgain50"
    "E1"  = "spkscpIsoStcDgain50Indisig"
    "F1"  = "spkfcpIsoStcDlat50"
    "G1"  = "spkfcpIsoStcDlat50Indisig"
    "H1"  = "spkfcpIsoStcDgain50"
    "I1"  = "spkfcpIsoStcDgain50Indisig"
    "J1"  = "spkscfIsoStcDlat50"
    "K1"  = "spkscfIsoStcDlat50Indisig"
    "L1"  = "spkscfIsoStcDgain50"
    "M1"  = "spkscfIsoStcDgain50Indisig"
    "N1"  = "spkfcfIsoStcDlat50"
    "O1"  = "spkfcfIsoStcDlat50Indisig"
    "P1"  = "spkfcfIsoStcDgain50"
    "Q1"  = "spkfcfIsoStcDgain50Indisig"
    "R1"  = "spkscpCrossStcDlat50"
    "S1"  = "spkscpCrossStcDlat50Indisig"
    "T1"  = "spkscpCrossStcDgain50"
    "U1"  = "spkscpCrossStcDgain50Indisig"
    "V1"  = "spkfcpCrossStcDlat50"
    "W1"  = "spkfcpCrossStcDlat50Indisig"
    "X1"  = "spkfcpCrossStcDgain50"
    "Y1"  = "spkfcpCrossStcDgain50Indisig"
    "Z1"  = "spksrndIsoStcDlat50"
    "AA1" = "spksrndIsoStcDlat50Indisig"
    "AB1" = "spksrndIsoStcDgain50"
    "AC1" = "spksrndIsoStcDgain50Indisig"
    "AD1" = "spkfrndIsoStcDlat50"
    "AE1" = "spkfrndIsoStcDlat50Indisig"
    "AF1" = "spkfrndIsoStcDgain50"
    "AG1" = "spkfrndIsoStcDgain50Indisig"
}

foreach ($addr in $renames.Keys) {
    $ws.Range($addr).Value = $renames[$addr]
}
